$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. The existing "Item Unit Name" column
# (with its header's themed font colour + fill) shifts from A to B, and a
# blank column A is created for the new "ID" header.
$ws.Columns.Item(1).Insert()

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Item Unit Name"

# Give the new "ID" header the same formatting (font colour/fill) as the
# "Item Unit Name" header by copying formats only.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Bold both header cells.
$ws.Range("A1:B1").Font.Bold = $true

# Restore / set the column widths (A narrower for "ID", B keeps the
# original 47-char width used for "Item Unit Name").
$ws.Columns.Item(1).ColumnWidth = 41
$ws.Columns.Item(2).ColumnWidth = 46.17

# Match the saved selection state (active cell B1).
$ws.Range("B1").Select()
